# Generate Report for Archive
#
# The localization status for e2e\fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md
# moved from "Ready for handoff" to "In Translation" on both the zh-cn and
# de-de targets. Update every sheet/table cell that surfaces that status.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: row 4 corresponds to fcee92e4-e213-4c85-bdcd-3f29e3f80fcf.md
# Columns E (zh-cn) and F (de-de) hold the per-locale status.
if ($overview.Range("E4").Text -eq $oldStatus) {
    $overview.Range("E4").Value = $newStatus
}
if ($overview.Range("F4").Text -eq $oldStatus) {
    $overview.Range("F4").Value = $newStatus
}

# zh-cn sheet: row 4 is the same file; column C is the Status column.
if ($zhcn.Range("C4").Text -eq $oldStatus) {
    $zhcn.Range("C4").Value = $newStatus
}

# de-de sheet: row 4 is the same file; column C is the Status column.
if ($dede.Range("C4").Text -eq $oldStatus) {
    $dede.Range("C4").Value = $newStatus
}

$wb.Save()
